# Costa Rica Primera Division - base update, 02-04-2024 23:59
# 1) Four pairs of rows that share an identical Date got re-sorted, swapping
#    everything except the running id in column A (same pattern Excel produces
#    when two same-timestamp rows trade places after a stable re-sort).
# 2) Two brand-new fixtures were appended as rows 222 and 223 (row 223 has no
#    result yet, so FTHG/FTAG/FTR/PL_AhOver/PL_AhUnder are left blank).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Swap row pairs that share the same Date (cols B:AC); column A (id) stays put ----
# rows 38 <-> 39
$ws.Range("B38").Value = 6781354
$ws.Range("C38").Value = "Costa Rica Primera Division"
$ws.Range("D38").Value = "Costa Rica Primera Division"
$ws.Range("E38").Value = 45171.75
$ws.Range("F38").Value = "Puntarenas"
$ws.Range("G38").Value = "AD San Carlos"
$ws.Range("H38").Value = 1
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = "H"
$ws.Range("K38").Value = 2.4
$ws.Range("L38").Value = 3.2
$ws.Range("M38").Value = 2.8
$ws.Range("N38").Value = 2.3
$ws.Range("O38").Value = 3.2
$ws.Range("P38").Value = 3
$ws.Range("Q38").Value = -0.25
$ws.Range("R38").Value = 2
$ws.Range("S38").Value = 1.8
$ws.Range("T38").Value = 2.25
$ws.Range("U38").Value = 1.9
$ws.Range("V38").Value = 1.9
$ws.Range("W38").Value = 1.3
$ws.Range("X38").Value = -1
$ws.Range("Y38").Value = -1
$ws.Range("Z38").Value = 1
$ws.Range("AA38").Value = -1
$ws.Range("AB38").Value = -1
$ws.Range("AC38").Value = 0.8999999999999999
$ws.Range("B39").Value = 6782522
$ws.Range("C39").Value = "Costa Rica Primera Division"
$ws.Range("D39").Value = "Costa Rica Primera Division"
$ws.Range("E39").Value = 45171.75
$ws.Range("F39").Value = "Municipal Perez Zeledon"
$ws.Range("G39").Value = "Sporting San Jose"
$ws.Range("H39").Value = 1
$ws.Range("I39").Value = 2
$ws.Range("J39").Value = "A"
$ws.Range("K39").Value = 2.5
$ws.Range("L39").Value = 3.5
$ws.Range("M39").Value = 2.5
$ws.Range("N39").Value = 2.2
$ws.Range("O39").Value = 3.5
$ws.Range("P39").Value = 2.9
$ws.Range("Q39").Value = -0.25
$ws.Range("R39").Value = 1.9
$ws.Range("S39").Value = 1.9
$ws.Range("T39").Value = 2.5
$ws.Range("U39").Value = 1.9
$ws.Range("V39").Value = 1.9
$ws.Range("W39").Value = -1
$ws.Range("X39").Value = -1
$ws.Range("Y39").Value = 1.9
$ws.Range("Z39").Value = -1
$ws.Range("AA39").Value = 0.8999999999999999
$ws.Range("AB39").Value = 0.8999999999999999
$ws.Range("AC39").Value = -1

# rows 110 <-> 111
$ws.Range("B110").Value = 6782581
$ws.Range("C110").Value = "Costa Rica Primera Division"
$ws.Range("D110").Value = "Costa Rica Primera Division"
$ws.Range("E110").Value = 45238.875
$ws.Range("F110").Value = "Alajuelense"
$ws.Range("G110").Value = "AD Grecia"
$ws.Range("H110").Value = 2
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = "H"
$ws.Range("K110").Value = 1.181
$ws.Range("L110").Value = 6.5
$ws.Range("M110").Value = 11
$ws.Range("N110").Value = 1.25
$ws.Range("O110").Value = 5
$ws.Range("P110").Value = 9
$ws.Range("Q110").Value = -1.75
$ws.Range("R110").Value = 1.975
$ws.Range("S110").Value = 1.825
$ws.Range("T110").Value = 3.25
$ws.Range("U110").Value = 2
$ws.Range("V110").Value = 1.8
$ws.Range("W110").Value = 0.25
$ws.Range("X110").Value = -1
$ws.Range("Y110").Value = -1
$ws.Range("Z110").Value = 0.4875
$ws.Range("AA110").Value = -0.5
$ws.Range("AB110").Value = -1
$ws.Range("AC110").Value = 0.8
$ws.Range("B111").Value = 6782579
$ws.Range("C111").Value = "Costa Rica Primera Division"
$ws.Range("D111").Value = "Costa Rica Primera Division"
$ws.Range("E111").Value = 45238.875
$ws.Range("F111").Value = "Santos de Gupiles"
$ws.Range("G111").Value = "AD San Carlos"
$ws.Range("H111").Value = 0
$ws.Range("I111").Value = 2
$ws.Range("J111").Value = "A"
$ws.Range("K111").Value = 2.4
$ws.Range("L111").Value = 3.3
$ws.Range("M111").Value = 2.7
$ws.Range("N111").Value = 2.375
$ws.Range("O111").Value = 3.4
$ws.Range("P111").Value = 2.8
$ws.Range("Q111").Value = -0.25
$ws.Range("R111").Value = 2
$ws.Range("S111").Value = 1.8
$ws.Range("T111").Value = 2.5
$ws.Range("U111").Value = 1.875
$ws.Range("V111").Value = 1.925
$ws.Range("W111").Value = -1
$ws.Range("X111").Value = -1
$ws.Range("Y111").Value = 1.8
$ws.Range("Z111").Value = -1
$ws.Range("AA111").Value = 0.8
$ws.Range("AB111").Value = -1
$ws.Range("AC111").Value = 0.925

# rows 129 <-> 130
$ws.Range("B129").Value = 6782598
$ws.Range("C129").Value = "Costa Rica Primera Division"
$ws.Range("D129").Value = "Costa Rica Primera Division"
$ws.Range("E129").Value = 45255.95833333334
$ws.Range("F129").Value = "Municipal Perez Zeledon"
$ws.Range("G129").Value = "Cartagines"
$ws.Range("H129").Value = 1
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = "H"
$ws.Range("K129").Value = 4.5
$ws.Range("L129").Value = 3.75
$ws.Range("M129").Value = 1.615
$ws.Range("N129").Value = 3.4
$ws.Range("O129").Value = 3.4
$ws.Range("P129").Value = 1.85
$ws.Range("Q129").Value = 0.5
$ws.Range("R129").Value = 1.8
$ws.Range("S129").Value = 2
$ws.Range("T129").Value = 2.75
$ws.Range("U129").Value = 1.9
$ws.Range("V129").Value = 1.9
$ws.Range("W129").Value = 2.4
$ws.Range("X129").Value = -1
$ws.Range("Y129").Value = -1
$ws.Range("Z129").Value = 0.8
$ws.Range("AA129").Value = -1
$ws.Range("AB129").Value = -1
$ws.Range("AC129").Value = 0.8999999999999999
$ws.Range("B130").Value = 6782596
$ws.Range("C130").Value = "Costa Rica Primera Division"
$ws.Range("D130").Value = "Costa Rica Primera Division"
$ws.Range("E130").Value = 45255.95833333334
$ws.Range("F130").Value = "Alajuelense"
$ws.Range("G130").Value = "AD Guanacasteca"
$ws.Range("H130").Value = 3
$ws.Range("I130").Value = 4
$ws.Range("J130").Value = "A"
$ws.Range("K130").Value = 1.363
$ws.Range("L130").Value = 4.75
$ws.Range("M130").Value = 8
$ws.Range("N130").Value = 1.444
$ws.Range("O130").Value = 4.333
$ws.Range("P130").Value = 7
$ws.Range("Q130").Value = -1.25
$ws.Range("R130").Value = 1.975
$ws.Range("S130").Value = 1.825
$ws.Range("T130").Value = 2.75
$ws.Range("U130").Value = 1.775
$ws.Range("V130").Value = 2.025
$ws.Range("W130").Value = -1
$ws.Range("X130").Value = -1
$ws.Range("Y130").Value = 6
$ws.Range("Z130").Value = -1
$ws.Range("AA130").Value = 0.825
$ws.Range("AB130").Value = 0.7749999999999999
$ws.Range("AC130").Value = -1

# rows 200 <-> 201
$ws.Range("B200").Value = 7624967
$ws.Range("C200").Value = "Costa Rica Primera Division"
$ws.Range("D200").Value = "Costa Rica Primera Division"
$ws.Range("E200").Value = 45353.75
$ws.Range("F200").Value = "Puntarenas"
$ws.Range("G200").Value = "Herediano"
$ws.Range("H200").Value = 0
$ws.Range("I200").Value = 0
$ws.Range("J200").Value = "D"
$ws.Range("K200").Value = 3.75
$ws.Range("L200").Value = 3.4
$ws.Range("M200").Value = 1.8
$ws.Range("N200").Value = 2.8
$ws.Range("O200").Value = 3.1
$ws.Range("P200").Value = 2.25
$ws.Range("Q200").Value = 0.25
$ws.Range("R200").Value = 1.8
$ws.Range("S200").Value = 2
$ws.Range("T200").Value = 2.25
$ws.Range("U200").Value = 1.775
$ws.Range("V200").Value = 2.025
$ws.Range("W200").Value = -1
$ws.Range("X200").Value = 2.1
$ws.Range("Y200").Value = -1
$ws.Range("Z200").Value = 0.4
$ws.Range("AA200").Value = -0.5
$ws.Range("AB200").Value = -1
$ws.Range("AC200").Value = 1.025
$ws.Range("B201").Value = 7623921
$ws.Range("C201").Value = "Costa Rica Primera Division"
$ws.Range("D201").Value = "Costa Rica Primera Division"
$ws.Range("E201").Value = 45353.75
$ws.Range("F201").Value = "AD Grecia"
$ws.Range("G201").Value = "Municipal Liberia"
$ws.Range("H201").Value = 1
$ws.Range("I201").Value = 2
$ws.Range("J201").Value = "A"
$ws.Range("K201").Value = 2.75
$ws.Range("L201").Value = 3.25
$ws.Range("M201").Value = 2.3
$ws.Range("N201").Value = 3.1
$ws.Range("O201").Value = 3.25
$ws.Range("P201").Value = 2.1
$ws.Range("Q201").Value = 0.25
$ws.Range("R201").Value = 1.9
$ws.Range("S201").Value = 1.9
$ws.Range("T201").Value = 2.5
$ws.Range("U201").Value = 1.9
$ws.Range("V201").Value = 1.9
$ws.Range("W201").Value = -1
$ws.Range("X201").Value = -1
$ws.Range("Y201").Value = 1.1
$ws.Range("Z201").Value = -1
$ws.Range("AA201").Value = 0.8999999999999999
$ws.Range("AB201").Value = 0.8999999999999999
$ws.Range("AC201").Value = -1

# ---- Append two brand-new rows (222, 223) ----
# Seed A222/A223 (id col, style 1) and E222/E223 (date col, style 2) by copying
# an existing same-styled cell, so the bold/border/center style and the date
# number format are reused instead of synthesizing new style entries.
$ws.Range("A221").Copy($ws.Range("A222"))
$ws.Range("A221").Copy($ws.Range("A223"))
$ws.Range("E221").Copy($ws.Range("E222"))
$ws.Range("E221").Copy($ws.Range("E223"))

# row 222
$ws.Range("A222").Value = 220
$ws.Range("B222").Value = 7623940
$ws.Range("C222").Value = "Costa Rica Primera Division"
$ws.Range("D222").Value = "Costa Rica Primera Division"
$ws.Range("E222").Value = 45383.75
$ws.Range("F222").Value = "AD Grecia"
$ws.Range("G222").Value = "Cartagines"
$ws.Range("H222").Value = 2
$ws.Range("I222").Value = 0
$ws.Range("J222").Value = "H"
$ws.Range("K222").Value = 3
$ws.Range("L222").Value = 3.4
$ws.Range("M222").Value = 2.1
$ws.Range("N222").Value = 3.1
$ws.Range("O222").Value = 3.3
$ws.Range("P222").Value = 2.1
$ws.Range("Q222").Value = 0.25
$ws.Range("R222").Value = 1.95
$ws.Range("S222").Value = 1.85
$ws.Range("T222").Value = 2.5
$ws.Range("U222").Value = 1.925
$ws.Range("V222").Value = 1.875
$ws.Range("W222").Value = 2.1
$ws.Range("X222").Value = -1
$ws.Range("Y222").Value = -1
$ws.Range("Z222").Value = 0.95
$ws.Range("AA222").Value = -1
$ws.Range("AB222").Value = -1
$ws.Range("AC222").Value = 0.875

# row 223
$ws.Range("A223").Value = 221
$ws.Range("B223").Value = 7980420
$ws.Range("C223").Value = "Costa Rica Primera Division"
$ws.Range("D223").Value = "Costa Rica Primera Division"
$ws.Range("E223").Value = 45385.95833333334
$ws.Range("F223").Value = "Deportivo Saprissa"
$ws.Range("G223").Value = "Municipal Perez Zeledon"
$ws.Range("K223").Value = 1.222
$ws.Range("L223").Value = 6
$ws.Range("M223").Value = 12
$ws.Range("N223").Value = 1.222
$ws.Range("O223").Value = 6
$ws.Range("P223").Value = 12
$ws.Range("Q223").Value = -1.75
$ws.Range("R223").Value = 1.825
$ws.Range("S223").Value = 1.975
$ws.Range("T223").Value = 3
$ws.Range("U223").Value = 1.9
$ws.Range("V223").Value = 1.9
$ws.Range("W223").Value = 0
$ws.Range("X223").Value = 0
$ws.Range("Y223").Value = 0
$ws.Range("Z223").Value = 0
$ws.Range("AA223").Value = 0

